$d = $word.ActiveDocument

# Update the date in the title paragraph
$d.Content.Find.Execute("2025-05-04 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-05 Monday", 2) | Out-Null

# Update the division problems in the table, identified by row/column
# since some values (e.g. "17÷4=") repeat with different replacements.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "88÷5="
$t.Cell(1, 2).Range.Text = "48÷5="
$t.Cell(1, 3).Range.Text = "54÷7="
$t.Cell(1, 4).Range.Text = "44÷6="
$t.Cell(1, 5).Range.Text = "36÷8="
$t.Cell(5, 1).Range.Text = "31÷6="
$t.Cell(5, 2).Range.Text = "58÷7="
$t.Cell(5, 3).Range.Text = "43÷9="
$t.Cell(5, 4).Range.Text = "82÷3="
$t.Cell(5, 5).Range.Text = "69÷3="
$t.Cell(9, 1).Range.Text = "39÷7="
$t.Cell(9, 2).Range.Text = "98÷6="
$t.Cell(9, 3).Range.Text = "68÷6="
$t.Cell(9, 4).Range.Text = "69÷9="
$t.Cell(9, 5).Range.Text = "37÷3="
$t.Cell(13, 1).Range.Text = "19÷2="
$t.Cell(13, 2).Range.Text = "85÷7="
$t.Cell(13, 3).Range.Text = "24÷8="
$t.Cell(13, 4).Range.Text = "42÷8="
$t.Cell(13, 5).Range.Text = "35÷2="
$t.Cell(17, 1).Range.Text = "10÷4="
$t.Cell(17, 2).Range.Text = "99÷9="
$t.Cell(17, 3).Range.Text = "26÷4="
$t.Cell(17, 4).Range.Text = "22÷8="
$t.Cell(17, 5).Range.Text = "76÷8="
